$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks and Bugs")

# Update Progress % values for the completed tasks (rows 17-19, 21, 23)
$ws.Range("E17").Value = 100
$ws.Range("E18").Value = 100
$ws.Range("E19").Value = 100
$ws.Range("E21").Value = 100
$ws.Range("E23").Value = 100

# Update the view: scroll position and current selection
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("E21").Select()
